# The paragraph "<id>p017r_1</id>" is currently split across three runs:
#   1) "<id>"      - Courier New, color 7f6000, sz/szCs 18
#   2) "p017r_1"   - default formatting, color 000000
#   3) "</id>"     - Courier New, color 7f6000, sz/szCs 18
# The edit merges them into a single run "<id>p017r_1</id>" that carries the
# Courier New / 7f6000 / 18 formatting of the former "<id>" and "</id>" runs.
#
# Rather than rewriting each run's font properties (which risks not
# reproducing every rPr attribute exactly), we delete the "p017r_1</id>"
# text and then re-insert it immediately after the remaining "<id>" run via
# InsertAfter. Word/this COM host merges appended text into the adjoining
# run when the run's own formatting applies, which reproduces the target
# run's rPr exactly (including attributes we would otherwise have to set
# one-by-one) and collapses the three runs into one.

$d = $word.ActiveDocument

$full = $d.Content.Text

$idStart = $full.IndexOf("<id>p017r_1</id>")
if ($idStart -lt 0) {
    throw "Could not find '<id>p017r_1</id>' in the document"
}

$idOpenEnd  = $idStart + 4                 # end of "<id>"
$idCloseEnd = $idStart + 16                # end of "<id>p017r_1</id>"

# Remove "p017r_1</id>" (the middle run plus the old closing-tag run).
$rRemove = $d.Range($idOpenEnd, $idCloseEnd)
$rRemove.Delete()

# Re-append the removed text to the remaining "<id>" run so it inherits
# that run's formatting and the three runs collapse into one.
$rOpen = $d.Range($idStart, $idOpenEnd)
$rOpen.InsertAfter("p017r_1</id>")
